# Updated cryptos list on Thu Mar  2 16:55:50 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Cells.Item(2, 4)
$cell.NumberFormat = '@'
$cell.Value = '23.316.73'
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(2, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -1.74%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(3, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.627.20'
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(3, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -1.85%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(4, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +0.21%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(5, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +0.14%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = '@'
$cell.Value = '297.52'
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(6, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -1.84%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(7, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.3748'
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(7, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -1.80%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(8, 4)
$cell.NumberFormat = '@'
$cell.Value = '50.39'
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(8, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -1.75%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.3463'
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(9, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -4.34%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.08006'
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(10, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -2.28%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(11, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -3.05%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(12, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +0.23%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = '@'
$cell.Value = '21.77'
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(13, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -3.64%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = '@'
$cell.Value = '6.272'
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(14, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -3.54%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = '@'
$cell.Value = '7.189'
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(15, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -3.12%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.00001182'
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(16, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -3.95%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.631.15'
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(17, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -1.47%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = '@'
$cell.Value = '94.34'
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(18, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -3.49%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.06934'
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(19, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -1.00%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = '@'
$cell.Value = '6.581'
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(20, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -3.75%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = '@'
$cell.Value = '17.20'
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(21, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -2.69%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(22, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +0.10%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(23, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -4.04%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = '@'
$cell.Value = '23.343.54'
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(24, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -1.62%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = '@'
$cell.Value = '2.436'
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(25, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -3.09%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = '@'
$cell.Value = '2.996'
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(26, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -1.22%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(27, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -2.53%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = '@'
$cell.Value = '151.03'
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(28, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -1.01%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = '@'
$cell.Value = '5.155'
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(29, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -1.18%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = '@'
$cell.Value = '131.10'
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(30, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -2.38%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.815.96'
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(31, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -1.08%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = '@'
$cell.Value = '6.665'
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(32, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -5.36%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(33, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -3.93%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = '@'
$cell.Value = '11.26'
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(34, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -5.93%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.9708'
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(35, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -8.56%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.02650'
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(36, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -5.75%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(37, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -1.11%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.2404'
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(38, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -4.61%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = '@'
$cell.Value = '5.812'
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(39, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -4.80%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.06641'
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(40, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -5.73%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = '@'
$cell.Value = '12.59'
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(41, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -3.56%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.6756'
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(42, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -3.77%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.291'
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(43, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -3.42%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = '@'
$cell.Value = '15.22'
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(44, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -5.09%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(45, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +0.17%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.6269'
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(46, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -3.95%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(47, 2)
$cell.NumberFormat = '@'
$cell.Value = 'NEARProtocol'
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(47, 3)
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = '@'
$cell.Value = '2.222'
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(47, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -4.01%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(48, 2)
$cell.NumberFormat = '@'
$cell.Value = 'PancakeSwap'
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(48, 3)
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = '@'
$cell.Value = '3.883'
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(48, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -2.04%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = '@'
$cell.Value = '126.33'
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(49, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -1.42%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.07617'
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(50, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -3.76%  '
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.210'
$cell.Style = 'Normal'

$cell = $ws.Cells.Item(51, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +0.71%  '
$cell.Style = 'Normal'
